$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8657803535461426
$ws.Range("C2").Value = 0.05122542381286621
$ws.Range("B3").Value = 0.7299034595489502
$ws.Range("C3").Value = 0.02620148658752441
$ws.Range("B4").Value = 0.70076584815979
$ws.Range("C4").Value = 0.02925825119018555
$ws.Range("B5").Value = 0.7345669269561768
$ws.Range("C5").Value = 0.03012681007385254
$ws.Range("B6").Value = 0.8423118591308594
$ws.Range("C6").Value = 0.04671835899353027
$ws.Range("B7").Value = 0.2211968898773193
$ws.Range("C7").Value = 0.007737874984741211
$ws.Range("B8").Value = 0.2312071323394775
$ws.Range("C8").Value = 0.006039857864379883
$ws.Range("B9").Value = 0.2614874839782715
$ws.Range("C9").Value = 0.005982637405395508
$ws.Range("B10").Value = 0.2686727046966553
$ws.Range("C10").Value = 0.004987239837646484
$ws.Range("B11").Value = 0.2399487495422363
$ws.Range("C11").Value = 0.007350444793701172
$ws.Range("B12").Value = 0.5863428115844727
$ws.Range("C12").Value = 0.04588222503662109
$ws.Range("B13").Value = 0.5235490798950195
$ws.Range("C13").Value = 0.04365944862365723
$ws.Range("B14").Value = 0.50528883934021
$ws.Range("C14").Value = 0.04187989234924316
$ws.Range("B15").Value = 0.5484380722045898
$ws.Range("C15").Value = 0.04631161689758301
$ws.Range("B16").Value = 0.603208065032959
$ws.Range("C16").Value = 0.04237270355224609
$ws.Range("B17").Value = 2.458347082138062
$ws.Range("C17").Value = 0.05089139938354492
$ws.Range("B18").Value = 2.406069993972778
$ws.Range("C18").Value = 0.02800774574279785
$ws.Range("B19").Value = 2.095066547393799
$ws.Range("C19").Value = 0.02736568450927734
$ws.Range("B20").Value = 2.166279792785645
$ws.Range("C20").Value = 0.02931594848632812
$ws.Range("B21").Value = 2.101835489273071
$ws.Range("C21").Value = 0.02836465835571289
$ws.Range("B22").Value = 0.3002429008483887
$ws.Range("C22").Value = 0.005982637405395508
$ws.Range("B23").Value = 0.3130626678466797
$ws.Range("C23").Value = 0.005982875823974609
$ws.Range("B24").Value = 0.3165256977081299
$ws.Range("C24").Value = 0.004971981048583984
$ws.Range("B25").Value = 0.3273334503173828
$ws.Range("C25").Value = 0.005983352661132812
$ws.Range("B26").Value = 0.3226416110992432
$ws.Range("C26").Value = 0.005247831344604492
$ws.Range("B27").Value = 0.8427574634552002
$ws.Range("C27").Value = 0.04801821708679199
$ws.Range("B28").Value = 0.9223864078521729
$ws.Range("C28").Value = 0.04850554466247559
$ws.Range("B29").Value = 0.8076572418212891
$ws.Range("C29").Value = 0.04718184471130371
$ws.Range("B30").Value = 0.8806717395782471
$ws.Range("C30").Value = 0.0488746166229248
$ws.Range("B31").Value = 0.843618631362915
$ws.Range("C31").Value = 0.04787135124206543
$ws.Range("B32").Value = 2.938529491424561
$ws.Range("C32").Value = 0.02655601501464844
$ws.Range("B33").Value = 2.943948030471802
$ws.Range("C33").Value = 0.02695322036743164
$ws.Range("B34").Value = 3.116305828094482
$ws.Range("C34").Value = 0.02593159675598145
$ws.Range("B35").Value = 2.972084045410156
$ws.Range("C35").Value = 0.024932861328125
$ws.Range("B36").Value = 2.960963726043701
$ws.Range("C36").Value = 0.02692818641662598
$ws.Range("B37").Value = 0.3869633674621582
$ws.Range("C37").Value = 0.004987955093383789
$ws.Range("B38").Value = 0.3979356288909912
$ws.Range("C38").Value = 0.005984783172607422
$ws.Range("B39").Value = 0.3901298046112061
$ws.Range("C39").Value = 0.00600123405456543
$ws.Range("B40").Value = 0.3948736190795898
$ws.Range("C40").Value = 0.005010128021240234
$ws.Range("B41").Value = 0.3976535797119141
$ws.Range("C41").Value = 0.005981922149658203
$ws.Range("B42").Value = 0.895395040512085
$ws.Range("C42").Value = 0.04595398902893066
$ws.Range("B43").Value = 0.8732194900512695
$ws.Range("C43").Value = 0.04443144798278809
$ws.Range("B44").Value = 0.8133976459503174
$ws.Range("C44").Value = 0.0438995361328125
$ws.Range("B45").Value = 0.8989794254302979
$ws.Range("C45").Value = 0.04589366912841797
$ws.Range("B46").Value = 0.8760039806365967
$ws.Range("C46").Value = 0.04485082626342773
$ws.Range("B47").Value = 3.658556938171387
$ws.Range("C47").Value = 0.02590513229370117
$ws.Range("B48").Value = 3.548771858215332
$ws.Range("C48").Value = 0.0269775390625
$ws.Range("B49").Value = 3.643791198730469
$ws.Range("C49").Value = 0.0259251594543457
$ws.Range("B50").Value = 3.587129831314087
$ws.Range("C50").Value = 0.02589941024780273
$ws.Range("B51").Value = 3.58620023727417
$ws.Range("C51").Value = 0.04695987701416016
$ws.Range("B52").Value = 0.7154374122619629
$ws.Range("C52").Value = 0.006982326507568359
$ws.Range("B53").Value = 0.7470018863677979
$ws.Range("C53").Value = 0.005990505218505859
$ws.Range("B54").Value = 0.5455429553985596
$ws.Range("C54").Value = 0.005982637405395508
$ws.Range("B55").Value = 0.6482686996459961
$ws.Range("C55").Value = 0.005982875823974609
$ws.Range("B56").Value = 0.6023893356323242
$ws.Range("C56").Value = 0.005984067916870117
$ws.Range("B57").Value = 1.627213001251221
$ws.Range("C57").Value = 0.04787087440490723
$ws.Range("B58").Value = 1.524784564971924
$ws.Range("C58").Value = 0.04690146446228027
$ws.Range("B59").Value = 1.517983436584473
$ws.Range("C59").Value = 0.04886889457702637
$ws.Range("B60").Value = 1.811950445175171
$ws.Range("C60").Value = 0.05089426040649414
$ws.Range("B61").Value = 1.524324178695679
$ws.Range("C61").Value = 0.04997515678405762

Write-Output "applied updates"
